$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$style = $wb.Styles.Add("Normal_GDP1")
$fnt = $style.Font
$fnt.Name = "Arial Cyr"
$fnt.Size = 10
$fnt.ColorIndex = -4105
$ws.Range("C6").Style = "Normal_GDP1"
